$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.119259259259259
$ws.Range("C2").Value = 0.0407407407407407
$ws.Range("D2").Value = 0.0155555555555556
$ws.Range("E2").Value = 0.00592592592592593
$ws.Range("F2").Value = 0.02
$ws.Range("G2").Value = 0.0103703703703704
$ws.Range("H2").Value = 0.00444444444444444
$ws.Range("I2").Value = 0.00444444444444444
$ws.Range("J2").Value = 0.017037037037037
$ws.Range("K2").Value = 0.00222222222222222
$ws.Range("L2").Value = 0.0111111111111111
$ws.Range("M2").Value = 0.979259259259259
$ws.Range("N2").Value = 0.00592592592592593
$ws.Range("O2").Value = 0.0407407407407407
$ws.Range("P2").Value = 0.0244444444444444
$ws.Range("Q2").Value = 0.00666666666666667
$ws.Range("R2").Value = 0.00222222222222222
$ws.Range("S2").Value = 0.0037037037037037
$ws.Range("T2").Value = 0.995555555555556
$ws.Range("U2").Value = 0.00592592592592593
$ws.Range("V2").Value = 0.00148148148148148
$ws.Range("W2").Value = 0.00888888888888889
$ws.Range("X2").Value = 0.000740740740740741
$ws.Range("B3").Value = 0.0177777777777778
$ws.Range("C3").Value = 0.900740740740741
$ws.Range("D3").Value = 0.951111111111111
$ws.Range("E3").Value = 0.000740740740740741
$ws.Range("F3").Value = 0.00666666666666667
$ws.Range("G3").Value = 0.0185185185185185
$ws.Range("H3").Value = 0.0111111111111111
$ws.Range("I3").Value = 0.952592592592593
$ws.Range("J3").Value = 0.00666666666666667
$ws.Range("K3").Value = 0.00666666666666667
$ws.Range("L3").Value = 0.00740740740740741
$ws.Range("N3").Value = 0.964444444444444
$ws.Range("O3").Value = 0.0851851851851852
$ws.Range("R3").Value = 0.0222222222222222
$ws.Range("U3").Value = 0.00666666666666667
$ws.Range("V3").Value = 0.0244444444444444
$ws.Range("W3").Value = 0.0266666666666667
$ws.Range("X3").Value = 0.0251851851851852
$ws.Range("B4").Value = 0.855555555555556
$ws.Range("C4").Value = 0.0214814814814815
$ws.Range("D4").Value = 0.0266666666666667
$ws.Range("E4").Value = 0.0466666666666667
$ws.Range("F4").Value = 0.971851851851852
$ws.Range("G4").Value = 0.964444444444444
$ws.Range("H4").Value = 0.975555555555556
$ws.Range("I4").Value = 0.0155555555555556
$ws.Range("J4").Value = 0.967407407407407
$ws.Range("K4").Value = 0.985925925925926
$ws.Range("L4").Value = 0.980740740740741
$ws.Range("M4").Value = 0.0192592592592593
$ws.Range("N4").Value = 0.0207407407407407
$ws.Range("O4").Value = 0.0192592592592593
$ws.Range("P4").Value = 0.974074074074074
$ws.Range("Q4").Value = 0.993333333333333
$ws.Range("R4").Value = 0.974814814814815
$ws.Range("S4").Value = 0.995555555555556
$ws.Range("T4").Value = 0.00222222222222222
$ws.Range("U4").Value = 0.986666666666667
$ws.Range("V4").Value = 0.0222222222222222
$ws.Range("W4").Value = 0.956296296296296
$ws.Range("X4").Value = 0.965185185185185
$ws.Range("B5").Value = 0.00666666666666667
$ws.Range("C5").Value = 0.0362962962962963
$ws.Range("D5").Value = 0.00666666666666667
$ws.Range("E5").Value = 0.946666666666667
$ws.Range("F5").Value = 0.000740740740740741
$ws.Range("G5").Value = 0.00666666666666667
$ws.Range("H5").Value = 0.00740740740740741
$ws.Range("I5").Value = 0.0274074074074074
$ws.Range("J5").Value = 0.00814814814814815
$ws.Range("K5").Value = 0.00518518518518519
$ws.Range("L5").Value = 0.000740740740740741
$ws.Range("M5").Value = 0.00148148148148148
$ws.Range("N5").Value = 0.00888888888888889
$ws.Range("O5").Value = 0.854814814814815
$ws.Range("P5").Value = 0.00148148148148148
$ws.Range("R5").Value = 0.000740740740740741
$ws.Range("S5").Value = 0.000740740740740741
$ws.Range("T5").Value = 0.00222222222222222
$ws.Range("U5").Value = 0.000740740740740741
$ws.Range("V5").Value = 0.951851851851852
$ws.Range("W5").Value = 0.00814814814814815
$ws.Range("X5").Value = 0.00814814814814815

Write-Output "Applied updated frequency table values"
